# Apply weekly update: insert 4 new daily-price rows for Tomate (Vega Monumental
# Concepción) just above the existing row 545, pushing the rest of the table
# (previously rows 545:610) down to 549:614.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows at the top of the block that starts at row 545.
$ws.Rows("545:548").Insert()

# Common/fixed values shared by every row of this sub-table.
$mercadoId = 11
$mercado   = "Vega Monumental Concepción"
$region    = "Bíobío"
$codreg    = 8
$catId     = 100112020
$categoria = "Tomate"
$clasif    = "Hortaliza"

# Row 545: Tomate, Larga vida, Primera - Provincia de Quillota
$ws.Cells.Item(545, 1).Value  = $mercadoId
$ws.Cells.Item(545, 2).Value  = $mercado
$ws.Cells.Item(545, 3).Value  = $region
$ws.Cells.Item(545, 4).Value  = 44946
$ws.Cells.Item(545, 5).Value  = $codreg
$ws.Cells.Item(545, 6).Value  = $catId
$ws.Cells.Item(545, 7).Value  = $categoria
$ws.Cells.Item(545, 8).Value  = "Larga vida"
$ws.Cells.Item(545, 9).Value  = "Primera"
$ws.Cells.Item(545, 10).Value = 150
$ws.Cells.Item(545, 11).Value = 14000
$ws.Cells.Item(545, 12).Value = 14000
$ws.Cells.Item(545, 13).Value = 14000
$ws.Cells.Item(545, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(545, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(545, 16).Value = 778
$ws.Cells.Item(545, 17).Value = 18
$ws.Cells.Item(545, 18).Value = $clasif

# Row 546: Tomate, Larga vida, Segunda - Provincia de Quillota
$ws.Cells.Item(546, 1).Value  = $mercadoId
$ws.Cells.Item(546, 2).Value  = $mercado
$ws.Cells.Item(546, 3).Value  = $region
$ws.Cells.Item(546, 4).Value  = 44946
$ws.Cells.Item(546, 5).Value  = $codreg
$ws.Cells.Item(546, 6).Value  = $catId
$ws.Cells.Item(546, 7).Value  = $categoria
$ws.Cells.Item(546, 8).Value  = "Larga vida"
$ws.Cells.Item(546, 9).Value  = "Segunda"
$ws.Cells.Item(546, 10).Value = 160
$ws.Cells.Item(546, 11).Value = 12000
$ws.Cells.Item(546, 12).Value = 12000
$ws.Cells.Item(546, 13).Value = 12000
$ws.Cells.Item(546, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(546, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(546, 16).Value = 667
$ws.Cells.Item(546, 17).Value = 18
$ws.Cells.Item(546, 18).Value = $clasif

# Row 547: Tomate, Semiduro, Primera - Quillón
$ws.Cells.Item(547, 1).Value  = $mercadoId
$ws.Cells.Item(547, 2).Value  = $mercado
$ws.Cells.Item(547, 3).Value  = $region
$ws.Cells.Item(547, 4).Value  = 44946
$ws.Cells.Item(547, 5).Value  = $codreg
$ws.Cells.Item(547, 6).Value  = $catId
$ws.Cells.Item(547, 7).Value  = $categoria
$ws.Cells.Item(547, 8).Value  = "Semiduro"
$ws.Cells.Item(547, 9).Value  = "Primera"
$ws.Cells.Item(547, 10).Value = 120
$ws.Cells.Item(547, 11).Value = 5000
$ws.Cells.Item(547, 12).Value = 5000
$ws.Cells.Item(547, 13).Value = 5000
$ws.Cells.Item(547, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(547, 15).Value = "Quillón"
$ws.Cells.Item(547, 16).Value = 500
$ws.Cells.Item(547, 17).Value = 10
$ws.Cells.Item(547, 18).Value = $clasif

# Row 548: Tomate, Semiduro, Segunda - Quillón
$ws.Cells.Item(548, 1).Value  = $mercadoId
$ws.Cells.Item(548, 2).Value  = $mercado
$ws.Cells.Item(548, 3).Value  = $region
$ws.Cells.Item(548, 4).Value  = 44946
$ws.Cells.Item(548, 5).Value  = $codreg
$ws.Cells.Item(548, 6).Value  = $catId
$ws.Cells.Item(548, 7).Value  = $categoria
$ws.Cells.Item(548, 8).Value  = "Semiduro"
$ws.Cells.Item(548, 9).Value  = "Segunda"
$ws.Cells.Item(548, 10).Value = 100
$ws.Cells.Item(548, 11).Value = 4000
$ws.Cells.Item(548, 12).Value = 4000
$ws.Cells.Item(548, 13).Value = 4000
$ws.Cells.Item(548, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(548, 15).Value = "Quillón"
$ws.Cells.Item(548, 16).Value = 400
$ws.Cells.Item(548, 17).Value = 10
$ws.Cells.Item(548, 18).Value = $clasif

# Give the new date cells (column D) the same date number format used
# throughout the rest of the column.
$ws.Range("D545:D548").NumberFormat = $ws.Range("D544").NumberFormat
